# Daily attendance processing - 2026-01-22 10:42:44
# For every "Recorded By" entry (column G) that lists multiple
# recorders/systems separated by ", ", rotate the list so the first
# entry moves to the end (e.g. "A, B" -> "B, A", "A, B, C" -> "B, C, A").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value2 = $rotated -join ", "
        }
    }
}
